# "#5: cash & deposit done"
# Adds the common trailing columns (property_category, category, date,
# legislator_name, legislator_id, source_file, index) to the 現金 (cash)
# and 存款 (deposit) sheets, matching the layout already used by the
# other sheets (土地/建物/汽車/股票/債務).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 4: 現金 (cash)
# ---------------------------------------------------------------------
$cash = $wb.Worksheets.Item(4)

# Header row
$cash.Cells.Item(1, 5).Value  = "property_category"
$cash.Cells.Item(1, 6).Value  = "category"
$cash.Cells.Item(1, 7).Value  = "date"
$cash.Cells.Item(1, 8).Value  = "legislator_name"
$cash.Cells.Item(1, 9).Value  = "legislator_id"
$cash.Cells.Item(1, 10).Value = "source_file"
$cash.Cells.Item(1, 11).Value = "index"

# Row 2 (index 40)
$cash.Cells.Item(2, 5).Value  = "cash"
$cash.Cells.Item(2, 6).Value  = "normal"
$cash.Cells.Item(2, 7).NumberFormat = "@"
$cash.Cells.Item(2, 7).Value  = "2013-12-24"
$cash.Cells.Item(2, 7).ClearFormats()
$cash.Cells.Item(2, 8).Value  = "林國正"
$cash.Cells.Item(2, 9).Value  = 1742
$cash.Cells.Item(2, 10).Value = "tmp399c1"
$cash.Cells.Item(2, 11).Value = 40

# Row 3 (index 41)
$cash.Cells.Item(3, 5).Value  = "cash"
$cash.Cells.Item(3, 6).Value  = "normal"
$cash.Cells.Item(3, 7).NumberFormat = "@"
$cash.Cells.Item(3, 7).Value  = "2013-12-24"
$cash.Cells.Item(3, 7).ClearFormats()
$cash.Cells.Item(3, 8).Value  = "林國正"
$cash.Cells.Item(3, 9).Value  = 1742
$cash.Cells.Item(3, 10).Value = "tmp399c1"
$cash.Cells.Item(3, 11).Value = 41

# ---------------------------------------------------------------------
# Sheet 5: 存款 (deposit)
# ---------------------------------------------------------------------
$deposit = $wb.Worksheets.Item(5)

# Header row
$deposit.Cells.Item(1, 7).Value  = "property_category"
$deposit.Cells.Item(1, 8).Value  = "category"
$deposit.Cells.Item(1, 9).Value  = "date"
$deposit.Cells.Item(1, 10).Value = "legislator_name"
$deposit.Cells.Item(1, 11).Value = "legislator_id"
$deposit.Cells.Item(1, 12).Value = "source_file"
$deposit.Cells.Item(1, 13).Value = "index"

# Data rows 2..12 correspond to index 46..56
$depositIndexes = @(46, 47, 48, 49, 50, 51, 52, 53, 54, 55, 56)
for ($i = 0; $i -lt $depositIndexes.Length; $i++) {
    $row = $i + 2
    $idx = $depositIndexes[$i]

    $deposit.Cells.Item($row, 7).Value  = "deposit"
    $deposit.Cells.Item($row, 8).Value  = "normal"
    $deposit.Cells.Item($row, 9).NumberFormat = "@"
    $deposit.Cells.Item($row, 9).Value  = "2013-12-24"
    $deposit.Cells.Item($row, 9).ClearFormats()
    $deposit.Cells.Item($row, 10).Value = "林國正"
    $deposit.Cells.Item($row, 11).Value = 1742
    $deposit.Cells.Item($row, 12).Value = "tmp399c1"
    $deposit.Cells.Item($row, 13).Value = $idx
}
